$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect Price column (D) cells from Excel auto-converting numeric-looking
# text (e.g. "1.00", "558.80") into actual numbers, which would silently
# drop significant trailing zeros / change the stored representation.
# These cells hold plain text in the source data (prices use "." as a
# thousands separator in some rows, e.g. "65.552.36"), so they must stay text.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '65.552.36'
$ws.Range('E2').Value = '  +1.26%  '
$ws.Range('D3').Value = '3.397.11'
$ws.Range('E3').Value = '  +0.86%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '558.80'
$ws.Range('E5').Value = '  +0.73%  '
$ws.Range('D6').Value = '175.62'
$ws.Range('E6').Value = '  -0.18%  '
$ws.Range('E7').Value = '  +1.81%  '
$ws.Range('D8').Value = '3.388.85'
$ws.Range('E8').Value = '  +0.86%  '
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').Value = '0.173'
$ws.Range('E10').Value = '  +5.17%  '
$ws.Range('D11').Value = '0.638'
$ws.Range('E11').Value = '  +1.08%  '
$ws.Range('D12').Value = '53.95'
$ws.Range('E12').Value = '  -1.65%  '
$ws.Range('D13').Value = '0.0000281'
$ws.Range('E13').Value = '  +1.94%  '
$ws.Range('D14').Value = '9.21'
$ws.Range('E14').Value = '  +1.21%  '
$ws.Range('D15').Value = '3.923.86'
$ws.Range('E15').Value = '  +0.64%  '
$ws.Range('D16').Value = '18.37'
$ws.Range('E16').Value = '  -0.13%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.398.84'
$ws.Range('E17').Value = '  +1.13%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = '0.119'
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').Value = '65.343.37'
$ws.Range('E19').Value = '  +1.17%  '
$ws.Range('D20').Value = '11.89'
$ws.Range('E20').Value = '  +0.01%  '
$ws.Range('D21').Value = '1.00'
$ws.Range('E21').Value = '  +1.93%  '
$ws.Range('D22').Value = '465.14'
$ws.Range('E22').Value = '  -0.11%  '
$ws.Range('E23').Value = '  +1.88%  '
$ws.Range('D24').Value = '4.10'
$ws.Range('E24').Value = '  +0.25%  '
$ws.Range('D25').Value = '14.12'
$ws.Range('E25').Value = '  +5.54%  '
$ws.Range('D26').Value = '87.91'
$ws.Range('E26').Value = '  +1.76%  '
$ws.Range('D27').Value = '2.91'
$ws.Range('E27').Value = '  +2.22%  '
$ws.Range('D28').Value = '10.75'
$ws.Range('E28').Value = '  -1.43%  '
$ws.Range('D29').Value = '8.76'
$ws.Range('E29').Value = '  -0.14%  '
$ws.Range('D30').Value = '31.17'
$ws.Range('E30').Value = '  +3.14%  '
$ws.Range('D31').Value = '6.58'
$ws.Range('E31').Value = '  -1.52%  '
$ws.Range('D32').Value = '63.47'
$ws.Range('E32').Value = '  +7.20%  '
$ws.Range('D33').Value = '11.51'
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D34').Value = '581.46'
$ws.Range('E34').Value = '  -0.60%  '
$ws.Range('E35').Value = '  -0.23%  '
$ws.Range('E36').Value = '  +0.04%  '
$ws.Range('D37').Value = '3.65'
$ws.Range('E37').Value = '  +4.65%  '
$ws.Range('D38').Value = '0.143'
$ws.Range('E38').Value = '  +1.79%  '
$ws.Range('D39').Value = '35.87'
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('D40').Value = '0.374'
$ws.Range('E40').Value = '  +0.04%  '
$ws.Range('D41').Value = '0.0₃0745'
$ws.Range('E41').Value = '  -1.99%  '
$ws.Range('D42').Value = '3.100.26'
$ws.Range('E42').Value = '  -0.61%  '
$ws.Range('D43').Value = '0.0419'
$ws.Range('E43').Value = '  +1.54%  '
$ws.Range('D44').Value = '2.78'
$ws.Range('E44').Value = '  -0.98%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').Value = '2.46'
$ws.Range('E45').Value = '  -2.74%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').Value = '0.134'
$ws.Range('E46').Value = '  +2.33%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '3.16'
$ws.Range('E47').Value = '  -2.00%  '
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('D49').Value = '140.60'
$ws.Range('E49').Value = '  +3.48%  '
$ws.Range('D50').Value = '2.54'
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('D51').Value = '8.38'
$ws.Range('E51').Value = '  -0.22%  '
